$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The A:D crosswalk block (material/disposition/swims_material/swims_disp) for
# each material grouping shifts up by one row starting at row 13, so that the
# first entry of each group (previously duplicated across two rows) collapses
# into a single row. Row 42's old A:D values are dropped (no row to replace
# them), leaving the last row's crosswalk columns empty.
for ($r = 13; $r -le 42; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r - 1, $c).Value = $v
    }
}
$ws.Range("A42:D42").ClearContents()

# Update the view: scroll the frozen pane back to the top and select A12:D12.
$ws.Range("A12:D12").Select()
